# PWLPOS.v2 user.xlsx template update
#
# The "user" import template previously listed sample staff accounts
# (staff10 / Staff 10, staff20 / Staff 20, manager6 / Manager 6) all on
# level_id 3/2. The refreshed template instead lists 5 sample customer
# accounts (customer1..customer5) on level_id 4, paired with real-looking
# display names, and grows the sample table from 3 rows to 5 rows.
#
# Columns: A=level_id, B=username, C=nama, D=password

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "nama" (display name) column first for the rows that already
# existed (2-4), then the "username" column for every data row (2-6),
# then the two new rows' "nama" values - this mirrors how the source
# data was written and keeps the shared-string table order stable.
$ws.Range("C2").Value = "Djoko Susanto"
$ws.Range("C3").Value = "Bachtiar Karim"
$ws.Range("C4").Value = "Susilo Wonowidjojo"

$ws.Range("B2").Value = "customer1"
$ws.Range("B3").Value = "customer2"
$ws.Range("B4").Value = "customer3"
$ws.Range("B5").Value = "customer4"
$ws.Range("B6").Value = "customer5"

$ws.Range("C5").Value = "Erick Thohir"
$ws.Range("C6").Value = "Arini Subianto"

# level_id becomes 4 for every sample row (2-6)
$ws.Range("A2").Value = 4
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 4

# password stays the same sample value for every row, including the two
# newly added ones
$ws.Range("D2").Value = 12345
$ws.Range("D3").Value = 12345
$ws.Range("D4").Value = 12345
$ws.Range("D5").Value = 12345
$ws.Range("D6").Value = 12345

# Selection ends up on C4 in the saved file
$ws.Range("C4").Select()
